$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MasterExecutor")

# Rows whose Runmode (column E) flips from "No" to "Yes"
$rowsToEnable = @(3,4,5,6,7,8,9,10,11,12,13,14,15,17,18,21,22,23,25,28,29)
foreach ($r in $rowsToEnable) {
    $ws.Range("E$r").Value = "Yes"
}

# Previously filtered-out (hidden) rows become visible again
$rowsToUnhide = @(16,19,20,24,26,27)
foreach ($r in $rowsToUnhide) {
    $ws.Rows($r).Hidden = $false
}

# Drop the "Yes only" criterion on the Runmode column, keeping the AutoFilter dropdowns
$ws.Range("A1:F29").AutoFilter(5)

# Restore the view: scrolled to top-left A7, selection on E4:E29
$ws.Range("E4:E29").Select()
$ws.Application.ActiveWindow.ScrollRow = 7
